# Fruta / hortaliza, semanal
# Insert one new weekly observation row ahead of the existing row 68,
# pushing the existing data (old rows 68-149) down by one row
# (they become new rows 69-150), and populate the newly inserted row
# with the new observation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 68 - shifts rows 68:149 down to 69:150
$ws.Rows(68).Insert()

# Populate the newly inserted row 68 with the new observation
$ws.Range("A68").Value = 10
$ws.Range("B68").Value = "Vega Modelo de Temuco"
$ws.Range("C68").Value = "La Araucanía"
$ws.Range("D68").Value = 44781
$ws.Range("E68").Value = 9
$ws.Range("F68").Value = 100112031
$ws.Range("G68").Value = "Poroto verde"
$ws.Range("H68").Value = "Sin especificar"
$ws.Range("I68").Value = "Primera"
$ws.Range("J68").Value = 40
$ws.Range("K68").Value = 30000
$ws.Range("L68").Value = 30000
$ws.Range("M68").Value = 30000
$ws.Range("N68").Value = "`$/malla 25 kilos"
$ws.Range("O68").Value = "Región de Arica y Parinacota"
$ws.Range("P68").Value = 1200
$ws.Range("Q68").Value = 25
$ws.Range("R68").Value = "Hortaliza"
